# Update 2p0. Convention change to support multi-axle vehicles
#
# Adds two new vehicle sheets (Truck_Amandla, Trailer_Kumanzi) modelled on
# the existing Trailer_Thwala sheet, inserted so the final tab order is:
#   Sedan_HambaLG, Sedan_Hamba, Bus_Makhulu, Truck_Amandla, Trailer_Thwala, Trailer_Kumanzi

$wb = $excel.ActiveWorkbook

$trailer = $wb.Worksheets.Item("Trailer_Thwala")

# --- Create "Truck_Amandla" just BEFORE "Trailer_Thwala" -------------------
# Copy(Before:=trailer) inserts the new copy immediately before $trailer,
# taking over $trailer's old tab slot (the $trailer handle now resolves to
# that new copy, while the original sheet - still named "Trailer_Thwala" -
# slides one slot later).
$trailer.Copy($trailer)
$truck = $trailer
$truck.Name = "Truck_Amandla"

# Drop the trailing (unused) rows 9:10 so the sheet matches the smaller
# 8-row layout used for the new vehicle entries.
$truck.Rows("9:10").Delete()

# Update the CAD instance/class labels for this vehicle.
$truck.Range("H3").Value = "CAD_Truck_Amandla"
$truck.Range("H4").Value = "CAD_Truck_Amandla"

# Update the Color (x,y,z) values.
$truck.Range("F7").Value = 0.6
$truck.Range("G7").Value = 0.8
$truck.Range("H7").Value = 1

# Update Opacity.
$truck.Range("H8").Value = 1

[void]$truck.Range("G23").Select()

# --- Create "Trailer_Kumanzi" just AFTER "Trailer_Thwala" ------------------
# Re-fetch the (real, original) Trailer_Thwala sheet by name.
$trailer = $wb.Worksheets.Item("Trailer_Thwala")
$trailer.Copy($null, $trailer)
$kumanzi = $wb.Worksheets.Item($trailer.Index + 1)
$kumanzi.Name = "Trailer_Kumanzi"

# Drop the trailing (unused) rows 9:10.
$kumanzi.Rows("9:10").Delete()

# Update the CAD instance/class labels for this vehicle.
$kumanzi.Range("H3").Value = "CAD_Trailer_Kumanzi"
$kumanzi.Range("H4").Value = "CAD_Trailer_Kumanzi"

# Update the Color (x,y,z) values.
$kumanzi.Range("F7").Value = 1
$kumanzi.Range("G7").Value = 0.75
$kumanzi.Range("H7").Value = 0.055

# Update Opacity.
$kumanzi.Range("H8").Value = 0.5

# Trailer_Kumanzi becomes the active/selected sheet (mirrors the diff's
# bookView activeTab pointing at the last tab).
$kumanzi.Activate()
[void]$kumanzi.Range("H8").Select()
